# Updates recalculated Leve profit figures across multiple sheets.
# Columns: H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#          K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ

$wb = $excel.ActiveWorkbook

function Set-LeveRow($SheetName, $Row, $H, $I, $J, $K, $L, $M, $N) {
    $ws = $wb.Worksheets.Item($SheetName)
    $ws.Cells.Item($Row, 8).Value = $H
    $ws.Cells.Item($Row, 9).Value = $I
    $ws.Cells.Item($Row, 10).Value = $J
    $ws.Cells.Item($Row, 11).Value = $K
    $ws.Cells.Item($Row, 12).Value = $L
    $ws.Cells.Item($Row, 13).Value = $M
    $ws.Cells.Item($Row, 14).Value = $N
}

# ALC
Set-LeveRow "ALC" 70  2667599.5 6667570   952.55554  20002710   2857.66662  -20002440   -3397.66662
Set-LeveRow "ALC" 73  2667599.5 6667570   952.55554  20002710   2857.66662  -20001774   -4729.66662
Set-LeveRow "ALC" 100 1716.0454 1461      2084.4443  1461       2084.4443   -920        -3166.4443
Set-LeveRow "ALC" 129 52904.6   308.8     79202.5    926.4000000000001 237607.5 4073.6   -247607.5
Set-LeveRow "ALC" 132 31566138  45912372  4419.8     137737116  13259.4     -137734586  -18319.4

# ARM
Set-LeveRow "ARM" 32  18801.148 13526.724 36508.145  13526.724  36508.145   -13239.724  -37082.145
Set-LeveRow "ARM" 61  4849.4517 5150.1924 3285.6     5150.1924  3285.6      -4938.1924  -3709.6
Set-LeveRow "ARM" 74  985.5714  748.4     1578.5     748.4      1578.5      125.6       -3326.5
Set-LeveRow "ARM" 77  985.5714  748.4     1578.5     3742       7892.5      626         -16628.5
Set-LeveRow "ARM" 102 142859000 200001620 2500       200001620  2500        -199999998  -5744
Set-LeveRow "ARM" 122 38463560  40001980  3014       120005940  9042        -120003490  -13942
Set-LeveRow "ARM" 136 4849.4517 5150.1924 3285.6     15450.5772 9856.799999999999 -12900.5772 -14956.8

# CRP
Set-LeveRow "CRP" 58  2150010   3426373   5720.24    3426373    5720.24     -3426170    -6126.24
Set-LeveRow "CRP" 62  4044.2856 3603.3333 4375       3603.3333  4375        -2979.3333  -5623
Set-LeveRow "CRP" 65  4044.2856 3603.3333 4375       18016.6665 21875       -14896.6665 -28115
Set-LeveRow "CRP" 136 2150010   3426373   5720.24    10279119   17160.72    -10276569   -22260.72

# CUL
Set-LeveRow "CUL" 132 2856.3086 657.8570999999999 3315.6865 5920.7139 29841.1785 -3390.7139 -34901.17849999999
Set-LeveRow "CUL" 137 37956.656 74347.71000000001 3991.6667 223043.13 11975.0001 -217943.13 -22175.0001

# CUL row 141 - M141 did not exist before (cell was absent/blank); now created.
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(141, 8).Value = 7191.857
$ws.Cells.Item(141, 9).Value = 3000
$ws.Cells.Item(141, 10).Value = 7890.5
$ws.Cells.Item(141, 11).Value = 9000
$ws.Cells.Item(141, 12).Value = 23671.5
$ws.Cells.Item(141, 13).Value = -3820
$ws.Cells.Item(141, 14).Value = -34031.5

# GSM
Set-LeveRow "GSM" 113 1090.6666 1113.4445 1067.8889  1113.4445  1067.8889   1056.5555   -5407.8889

# LTW
Set-LeveRow "LTW" 61  2444      1764.8889 5500       1764.8889  5500        -1562.8889  -5904
Set-LeveRow "LTW" 113 2444      1764.8889 5500       1764.8889  5500        405.1111000000001 -9840

# LTW row 122 - J, L, N are unchanged in this row; only H, I, K, M change.
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 13165099
$ws.Cells.Item(122, 9).Value = 10917.637
$ws.Cells.Item(122, 11).Value = 32752.911
$ws.Cells.Item(122, 13).Value = -30302.911

# WVR
Set-LeveRow "WVR" 113 1173.421  1000.4    1235.2142  3001.2     3705.6426   -831.1999999999998 -8045.642599999999

# WVR row 122 - J, L, N are unchanged in this row; only H, I, K, M change.
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 968.4375
$ws.Cells.Item(122, 9).Value = 916.4167
$ws.Cells.Item(122, 11).Value = 2749.2501
$ws.Cells.Item(122, 13).Value = -299.2501000000002
